$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "29.628.19"
Set-TextCell 2 5 "  +2.68%  "

Set-TextCell 3 4 "1.991.94"
Set-TextCell 3 5 "  +6.15%  "

Set-TextCell 4 5 "  -0.15%  "

Set-TextCell 5 4 "328.87"
Set-TextCell 5 5 "  +1.32%  "

Set-TextCell 6 4 "0.9999"
Set-TextCell 6 5 "  -0.26%  "

Set-TextCell 7 4 "0.4678"
Set-TextCell 7 5 "  +1.47%  "

Set-TextCell 8 5 "  +2.13%  "

Set-TextCell 9 2 "Dogecoin"
Set-TextCell 9 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell 9 4 "0.08078"
Set-TextCell 9 5 "  +3.03%  "

Set-TextCell 10 2 "Polygon"
Set-TextCell 10 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 10 4 "1.001"
Set-TextCell 10 5 "  +1.90%  "

Set-TextCell 11 2 "Solana"
Set-TextCell 11 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 11 4 "22.95"
Set-TextCell 11 5 "  +5.67%  "

Set-TextCell 12 2 "WrappedEther"
Set-TextCell 12 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 12 4 "1.984.35"
Set-TextCell 12 5 "  +4.71%  "

Set-TextCell 13 2 "Chainlink"
Set-TextCell 13 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 13 4 "7.241"
Set-TextCell 13 5 "  +3.63%  "

Set-TextCell 14 2 "Polkadot"
Set-TextCell 14 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 14 4 "5.869"
Set-TextCell 14 5 "  +3.68%  "

Set-TextCell 15 2 "TRON"
Set-TextCell 15 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 15 4 "0.07131"
Set-TextCell 15 5 "  +2.29%  "

Set-TextCell 16 2 "Litecoin"
Set-TextCell 16 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 16 4 "88.85"
Set-TextCell 16 5 "  +0.77%  "

Set-TextCell 17 2 "BinanceUSD"
Set-TextCell 17 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell 17 4 "1.002"
Set-TextCell 17 5 "  -0.14%  "

Set-TextCell 18 2 "ShibaInu"
Set-TextCell 18 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 18 4 "0.00001003"
Set-TextCell 18 5 "  +0.81%  "

Set-TextCell 19 2 "Avalanche"
Set-TextCell 19 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 19 4 "17.43"
Set-TextCell 19 5 "  +3.04%  "

Set-TextCell 20 2 "Dai"
Set-TextCell 20 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 20 4 "1.000"
Set-TextCell 20 5 "  -0.20%  "

Set-TextCell 21 2 "WrappedBTC"
Set-TextCell 21 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 21 4 "29.609.39"
Set-TextCell 21 5 "  +2.60%  "

Set-TextCell 22 2 "Uniswap"
Set-TextCell 22 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 22 4 "5.561"
Set-TextCell 22 5 "  +5.74%  "

Set-TextCell 23 2 "Cosmos"
Set-TextCell 23 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 23 4 "11.25"
Set-TextCell 23 5 "  +2.46%  "

Set-TextCell 24 2 "Toncoin"
Set-TextCell 24 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 24 4 "2.109"
Set-TextCell 24 5 "  +0.19%  "

Set-TextCell 25 2 "Monero"
Set-TextCell 25 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 25 4 "157.78"
Set-TextCell 25 5 "  +1.11%  "

Set-TextCell 26 2 "EthereumClassic"
Set-TextCell 26 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 26 4 "19.69"
Set-TextCell 26 5 "  +2.05%  "

Set-TextCell 27 2 "InternetComputer(DFINITY)"
Set-TextCell 27 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 27 4 "5.984"
Set-TextCell 27 5 "  +0.85%  "

Set-TextCell 28 2 "BitcoinCash"
Set-TextCell 28 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 28 4 "120.46"
Set-TextCell 28 5 "  +2.45%  "

Set-TextCell 29 2 "LidoDAOToken"
Set-TextCell 29 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell 29 4 "1.945"
Set-TextCell 29 5 "  +2.17%  "

Set-TextCell 30 2 "Stellar"
Set-TextCell 30 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 30 4 "0.09447"
Set-TextCell 30 5 "  +1.06%  "

Set-TextCell 31 2 "ImmutableX"
Set-TextCell 31 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 31 4 "0.9160"
Set-TextCell 31 5 "  +1.79%  "

Set-TextCell 32 2 "Filecoin"
Set-TextCell 32 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 32 4 "5.284"
Set-TextCell 32 5 "  +0.40%  "

Set-TextCell 33 2 "ARBITRUM"
Set-TextCell 33 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 33 4 "1.355"
Set-TextCell 33 5 "  +2.98%  "

Set-TextCell 34 2 "HuobiToken"
Set-TextCell 34 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 34 4 "3.185"
Set-TextCell 34 5 "  -2.02%  "

Set-TextCell 35 2 "Hedera"
Set-TextCell 35 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 35 4 "0.05858"
Set-TextCell 35 5 "  +2.27%  "

Set-TextCell 36 2 "TrustWalletToken"
Set-TextCell 36 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 36 4 "1.178"
Set-TextCell 36 5 "  +0.70%  "

Set-TextCell 37 2 "PEPE"
Set-TextCell 37 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell 37 4 "0.000003395"
Set-TextCell 37 5 "  +77.25%  "

Set-TextCell 38 4 "0.02130"
Set-TextCell 38 5 "  +2.87%  "

Set-TextCell 39 2 "FraxShare"
Set-TextCell 39 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 39 4 "7.907"
Set-TextCell 39 5 "  +3.90%  "

Set-TextCell 40 2 "TheSandbox"
Set-TextCell 40 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell 40 4 "0.5793"
Set-TextCell 40 5 "  +2.54%  "

Set-TextCell 41 2 "Algorand"
Set-TextCell 41 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell 41 4 "0.1826"
Set-TextCell 41 5 "  +3.20%  "

Set-TextCell 42 2 "Aptos"
Set-TextCell 42 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell 42 4 "9.903"
Set-TextCell 42 5 "  +2.31%  "

Set-TextCell 43 2 "MXToken"
Set-TextCell 43 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 43 4 "2.821"
Set-TextCell 43 5 "  +11.14%  "

Set-TextCell 44 2 "EnergySwap"
Set-TextCell 44 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell 44 4 "12.12"
Set-TextCell 44 5 "  +1.42%  "

Set-TextCell 45 2 "Decentraland"
Set-TextCell 45 3 "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell 45 4 "0.5404"
Set-TextCell 45 5 "  +1.41%  "

Set-TextCell 46 2 "RenderToken"
Set-TextCell 46 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 46 4 "2.212"
Set-TextCell 46 5 "  -0.65%  "

Set-TextCell 47 2 "Cronos"
Set-TextCell 47 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell 47 4 "0.06976"
Set-TextCell 47 5 "  -0.89%  "

Set-TextCell 48 2 "NEARProtocol"
Set-TextCell 48 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 48 4 "1.873"
Set-TextCell 48 5 "  +1.88%  "

Set-TextCell 49 2 "Quant"
Set-TextCell 49 3 "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell 49 4 "114.12"
Set-TextCell 49 5 "  +1.57%  "

Set-TextCell 50 2 "WOONetwork"
Set-TextCell 50 3 "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextCell 50 4 "0.3089"
Set-TextCell 50 5 "  +8.22%  "

Set-TextCell 51 2 "Aave"
Set-TextCell 51 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 51 4 "73.57"
Set-TextCell 51 5 "  +3.96%  "
